$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codebook used to have 10 separate rows for per-trait max/min scores
# (maxscore_bfi_a/c/e/n/o, minscore_bfi_a/c/e/n/o). These are replaced by a
# single "maxscore" row and a single "minscore" row, since the calculation
# now produces one overall mean value rather than five per-trait sums.
# Delete the 8 extra rows (rows 30-37), which shifts every row below them
# up by 8.
$ws.Range("A30:A37").EntireRow.Delete()

# Rename the remaining two rows and clear their (already-empty) explanation
# columns so they match the rest of the "no extra info" rows in the sheet.
$ws.Range("A28").Value = "maxscore"
$ws.Range("B28:D28").ClearContents()

$ws.Range("A29").Value = "minscore"
$ws.Range("B29:D29").ClearContents()
